$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2632689
$ws.Range("I19").Value = 4386761.5
$ws.Range("J19").Value = 1580.25
$ws.Range("K19").Value = 4386761.5
$ws.Range("L19").Value = 1580.25
$ws.Range("M19").Value = -4386586.5
$ws.Range("N19").Value = -1930.25

$ws.Range("H123").Value = 40952.5
$ws.Range("J123").Value = 41270
$ws.Range("L123").Value = 41270
$ws.Range("N123").Value = -51070

$ws.Range("H138").Value = 2497.78
$ws.Range("I138").Value = 623.3333
$ws.Range("J138").Value = 2909.244
$ws.Range("K138").Value = 1869.9999
$ws.Range("L138").Value = 8727.732
$ws.Range("M138").Value = 3270.0001
$ws.Range("N138").Value = -19007.732

$ws.Range("H141").Value = 50015.855
$ws.Range("I141").Value = 60850.176
$ws.Range("J141").Value = 3970
$ws.Range("K141").Value = 182550.528
$ws.Range("L141").Value = 11910
$ws.Range("M141").Value = -177370.528
$ws.Range("N141").Value = -22270

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 10856.571
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 10856.571
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 10856.571
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -12340.571

$ws.Range("H76").Value = 11144
$ws.Range("J76").Value = 20288
$ws.Range("L76").Value = 20288
$ws.Range("N76").Value = -20964

$ws.Range("H79").Value = 11144
$ws.Range("J79").Value = 20288
$ws.Range("L79").Value = 20288
$ws.Range("N79").Value = -22628

$ws.Range("H80").Value = 39351.715
$ws.Range("J80").Value = 39351.715
$ws.Range("L80").Value = 39351.715
$ws.Range("N80").Value = -41347.715

$ws.Range("H83").Value = 39351.715
$ws.Range("J83").Value = 39351.715
$ws.Range("L83").Value = 118055.145
$ws.Range("N83").Value = -128039.145

$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 45000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 2680.1667
$ws.Range("I132").Value = 1095.1818
$ws.Range("J132").Value = 5170.857
$ws.Range("K132").Value = 3285.5454
$ws.Range("L132").Value = 15512.571
$ws.Range("M132").Value = -755.5454
$ws.Range("N132").Value = -20572.571

$ws.Range("H137").Value = 40703.332
$ws.Range("J137").Value = 40703.332
$ws.Range("L137").Value = 40703.332
$ws.Range("N137").Value = -50903.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -46694

$ws.Range("H60").Value = 19000
$ws.Range("J60").Value = 19000
$ws.Range("L60").Value = 19000
$ws.Range("N60").Value = -20198

$ws.Range("H107").Value = 1829.0646
$ws.Range("I107").Value = 1455.4762
$ws.Range("J107").Value = 2613.6
$ws.Range("K107").Value = 1455.4762
$ws.Range("L107").Value = 2613.6
$ws.Range("M107").Value = 464.5237999999999
$ws.Range("N107").Value = -6453.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5570.75
$ws.Range("I31").Value = 1028.7826
$ws.Range("J31").Value = 17178
$ws.Range("K31").Value = 1028.7826
$ws.Range("L31").Value = 17178
$ws.Range("M31").Value = -733.7826
$ws.Range("N31").Value = -17768

$ws.Range("H34").Value = 5570.75
$ws.Range("I34").Value = 1028.7826
$ws.Range("J34").Value = 17178
$ws.Range("K34").Value = 1028.7826
$ws.Range("L34").Value = 17178
$ws.Range("M34").Value = -826.7826
$ws.Range("N34").Value = -17582

$ws.Range("H52").Value = 34600
$ws.Range("J52").Value = 34600
$ws.Range("L52").Value = 34600
$ws.Range("N52").Value = -35188

$ws.Range("H68").Value = 53282.125
$ws.Range("J68").Value = 53282.125
$ws.Range("L68").Value = 53282.125
$ws.Range("N68").Value = -54780.125

$ws.Range("H71").Value = 53282.125
$ws.Range("J71").Value = 53282.125
$ws.Range("L71").Value = 159846.375
$ws.Range("N71").Value = -167334.375

$ws.Range("H109").Value = 34071
$ws.Range("J109").Value = 34071
$ws.Range("L109").Value = 34071
$ws.Range("N109").Value = -36151

$ws.Range("H134").Value = 5119.5483
$ws.Range("I134").Value = 5727.909
$ws.Range("K134").Value = 17183.727
$ws.Range("M134").Value = -14648.727

$ws.Range("H138").Value = 43018
$ws.Range("J138").Value = 43018
$ws.Range("L138").Value = 43018
$ws.Range("N138").Value = -53298

$ws.Range("H140").Value = 104864.445
$ws.Range("J140").Value = 104864.445
$ws.Range("L140").Value = 104864.445
$ws.Range("N140").Value = -115224.445

$ws.Range("H141").Value = 29683.334
$ws.Range("J141").Value = 29683.334
$ws.Range("L141").Value = 29683.334
$ws.Range("N141").Value = -40043.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.933334
$ws.Range("I2").Value = 78.333336
$ws.Range("J2").Value = 144.33333
$ws.Range("K2").Value = 470.000016
$ws.Range("L2").Value = 865.9999799999999
$ws.Range("M2").Value = -357.000016
$ws.Range("N2").Value = -1091.99998

$ws.Range("H19").Value = 4444
$ws.Range("J19").Value = 4444
$ws.Range("L19").Value = 13332
$ws.Range("N19").Value = -13680

$ws.Range("H137").Value = 2202.2727
$ws.Range("I137").Value = 857.5
$ws.Range("J137").Value = 2970.7144
$ws.Range("K137").Value = 2572.5
$ws.Range("L137").Value = 8912.143199999999
$ws.Range("M137").Value = 2527.5
$ws.Range("N137").Value = -19112.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15936.75
$ws.Range("J46").Value = 15936.75
$ws.Range("L46").Value = 15936.75
$ws.Range("N46").Value = -16248.75

$ws.Range("H132").Value = 6025.273
$ws.Range("I132").Value = 4796
$ws.Range("J132").Value = 7049.6665
$ws.Range("K132").Value = 14388
$ws.Range("L132").Value = 21148.9995
$ws.Range("M132").Value = -11858
$ws.Range("N132").Value = -26208.9995

$ws.Range("H137").Value = 40277.2
$ws.Range("J137").Value = 40277.2
$ws.Range("L137").Value = 40277.2
$ws.Range("N137").Value = -50477.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 18004.166
$ws.Range("J43").Value = 18004.166
$ws.Range("L43").Value = 18004.166
$ws.Range("N43").Value = -18390.166

$ws.Range("H110").Value = 39002
$ws.Range("J110").Value = 39002
$ws.Range("L110").Value = 39002
$ws.Range("N110").Value = -47182

$ws.Range("H136").Value = 4028.2
$ws.Range("I136").Value = 1278.1875
$ws.Range("J136").Value = 8917.111000000001
$ws.Range("K136").Value = 3834.5625
$ws.Range("L136").Value = 26751.333
$ws.Range("M136").Value = -1284.5625
$ws.Range("N136").Value = -31851.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 67734.336
$ws.Range("J46").Value = 67734.336
$ws.Range("L46").Value = 67734.336
$ws.Range("N46").Value = -68196.336

$ws.Range("H122").Value = 7552.778
$ws.Range("I122").Value = 5495.8335
$ws.Range("K122").Value = 16487.5005
$ws.Range("M122").Value = -14037.5005

$ws.Range("H134").Value = 67734.336
$ws.Range("J134").Value = 67734.336
$ws.Range("L134").Value = 203203.008
$ws.Range("N134").Value = -208273.008

$ws.Range("H136").Value = 6660.4116
$ws.Range("I136").Value = 6159.5
$ws.Range("J136").Value = 7578.75
$ws.Range("K136").Value = 18478.5
$ws.Range("L136").Value = 22736.25
$ws.Range("M136").Value = -15928.5
$ws.Range("N136").Value = -27836.25
